# "Changes in E and F excel files" - this workbook is the "F suite" file.
# The "Test Cases" sheet gains a new "Jira id" column (inserted after TCID,
# before Description) populated with the OPQA-20x ticket ids, the test
# descriptions lose their leading "To " wording, and the active selection /
# used range grow to account for the new column.

$wb = $excel.ActiveWorkbook

# --- Window geometry (cosmetic; best effort) ------------------------------
try {
    $win = $wb.Windows.Item(1)
    $win.Left   = 4230
    $win.Top    = 3960
    $win.Width  = 14400
    $win.Height = 10125
} catch {
}

# --- "Test Cases" sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("Test Cases")

# Insert a new column B ("Jira id"), pushing Description/Runmode/Results
# from B/C/D to C/D/E.
$ws.Columns("B").Insert()

$ws.Range("B1").Value2 = "Jira id"
$ws.Range("B2").Value2 = "OPQA-206"
$ws.Range("B3").Value2 = "OPQA-207"
$ws.Range("B4").Value2 = "OPQA-208"
$ws.Range("B5").Value2 = "OPQA-209"

# The description column (now column C) drops the leading "To " wording.
$ws.Range("C2").Value2 = "Verify that user receives a notification when he is followed by someone"
$ws.Range("C3").Value2 = "Verify that user receives a notification when someone he is following comments on an article"
$ws.Range("C4").Value2 = "Verify that user receives a notification when someone comments on an article contained in his watchlist"
$ws.Range("C5").Value2 = "Verify that user receives a notification if someone likes his comment"

# Runmode cell in the last row keeps the plain border-only look (no fill
# flag), matching its neighbours after the column shift.
$ws.Range("D5").Interior.Pattern = -4142

# Selection now sits on the (new) Runmode column.
$ws.Range("D2:D5").Select()
